$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "NO_MATCHING_RECORD" error-code row (row 2) that documents the
# match-failure case now handled by the new Try-Catch around the member
# search.
$ws.Range("A2").Value = "NO_MATCHING_RECORD"
$ws.Range("B2").Value = "No match found in MRRS Report. Sailor may not be locally gained in MRRS."

# The columns were re-sized (bestfit) to accommodate the longer text that was
# just entered. Column A grows from ~9.8 to 22 characters wide and column B
# grows from ~20.3 to ~63.7 characters wide.
$ws.Columns.Item(1).ColumnWidth = 21.1
$ws.Columns.Item(2).ColumnWidth = 62.8

# The active cell / selection moved to the newly-entered cell B2.
$ws.Range("B2").Select()
